# Auto-generated update of loading_percent values (case with 380 kV)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(4.918926664783237, 4.632121578177307, 10.80326495706875, 23.53882848165438, 3.622697183321949, 0, 0, 0, 13.42095396178596, 0, 15.32938006460495, 0, 21.00210877370274),
    @(4.749674754588987, 4.600420658625072, 10.75119476843116, 23.58216135031835, 3.625106964426211, 0, 0, 0, 12.82613224543138, 0, 14.99886812371665, 0, 21.10102617189822),
    @(4.644144176141034, 4.580868679788154, 10.72424419551092, 23.61837789450059, 3.626663397786804, 0, 0, 0, 12.44688858813911, 0, 14.79515126053738, 0, 21.16912609868557),
    @(4.600812301464723, 4.572883276019547, 10.71452950512096, 23.63553892435155, 3.627317038679304, 0, 0, 0, 12.28901148792713, 0, 14.71205578184049, 0, 21.19871737711881),
    @(4.593599687212905, 4.571556351914058, 10.71299307096631, 23.63853316263091, 3.627426747771496, 0, 0, 0, 12.26260105423343, 0, 14.69825654595316, 0, 21.20374176833704),
    @(4.643561005888587, 4.58076105246634, 10.72410804104096, 23.61859962637516, 3.626672134454969, 0, 0, 0, 12.44477261858557, 0, 14.79403076928016, 0, 21.16951774274654),
    @(4.860947103975855, 4.621211770492003, 10.78427313721115, 23.55176718608135, 3.623512171565334, 0, 0, 0, 13.21887710969062, 0, 15.21565095371606, 0, 21.03468003203171),
    @(5.271389556882198, 4.699658419989746, 10.94172985881145, 23.49750447914906, 3.617922031815943, 0, 0, 0, 14.61847034035821, 0, 16.03128845998065, 0, 20.82921036248376),
    @(5.559597191190305, 4.756515345664814, 11.08081171865106, 23.50507427341135, 3.614180539860946, 0, 0, 0, 15.56658189321085, 0, 16.61741912546716, 0, 20.71487107904104),
    @(5.687159479792768, 4.782160836239849, 11.14897195344798, 23.51890958102457, 3.612556930349512, 0, 0, 0, 15.97921297757906, 0, 16.87998785734488, 0, 20.67094584886107),
    @(5.73490781169266, 4.791836311045429, 11.17546637424317, 23.52564711352601, 3.611953319751161, 0, 0, 0, 16.13269627475997, 0, 16.97873698098317, 0, 20.65548641313719),
    @(5.724649805316048, 4.789754191681585, 11.1697302574141, 23.5241293817981, 3.612082820262691, 0, 0, 0, 16.09976537387399, 0, 16.95750128264905, 0, 20.65876349069488),
    @(5.691099171355913, 4.782957588339332, 11.15113808460983, 23.51943383607117, 3.612507046521605, 0, 0, 0, 15.99189612326818, 0, 16.8881261582566, 0, 20.66965039976395),
    @(5.670474634610139, 4.778789664816379, 11.13983825022296, 23.51675289870268, 3.612768356094047, 0, 0, 0, 15.9254598863915, 0, 16.84554061855834, 0, 20.67647216659677),
    @(5.551185220660429, 4.754834566365903, 11.07645406948022, 23.50437961073033, 3.614288219248662, 0, 0, 0, 15.53923207700058, 0, 16.60016954090074, 0, 20.71790540484485),
    @(5.477061903567471, 4.740079634362246, 11.03880928456146, 23.49945417121977, 3.615240646420997, 0, 0, 0, 15.29744587946204, 0, 16.44853127516686, 0, 20.74540293755877),
    @(5.434096921471642, 4.731572635293838, 11.01761869764697, 23.49759920391264, 3.615795841635987, 0, 0, 0, 15.15662375227124, 0, 16.36093564015414, 0, 20.76197946796187),
    @(5.419494354352264, 4.728688951564412, 11.01052379429649, 23.49713895916605, 3.615985091314231, 0, 0, 0, 15.10864557761016, 0, 16.33121547037828, 0, 20.76772230036703),
    @(5.484987094290834, 4.741652460636362, 11.042768988803, 23.49987723525186, 3.615138495115949, 0, 0, 0, 15.32336654452438, 0, 16.46471318949694, 0, 20.74239696976741),
    @(5.700969262897684, 4.78495492523827, 11.15658067015173, 23.52077234378723, 3.612382137094257, 0, 0, 0, 16.02365576343253, 0, 16.90852249641674, 0, 20.66642070063103),
    @(5.838862989602039, 4.813044229369414, 11.2349373507616, 23.54316268287199, 3.610646039975129, 0, 0, 0, 16.46515248181339, 0, 17.19457450505397, 0, 20.62361411251164),
    @(5.765579213415601, 4.798073226335869, 11.19276029255049, 23.53041257746754, 3.611566668840496, 0, 0, 0, 16.23102307020686, 0, 17.04229930412242, 0, 20.64583066056408),
    @(5.481405203705764, 4.740941460869917, 11.04097739741174, 23.49968292579795, 3.615184653941297, 0, 0, 0, 15.31165346167035, 0, 16.45739863687119, 0, 20.74375357680246),
    @(5.162455552130699, 4.678555562725144, 10.89496671653174, 23.50389081328594, 3.619369811960744, 0, 0, 0, 14.25342701575302, 0, 15.81250386157653, 0, 20.87841275529016)
)

$arr = New-Object "object[,]" 24,13
for ($i = 0; $i -lt 24; $i++) {
    for ($j = 0; $j -lt 13; $j++) {
        $arr[$i,$j] = $data[$i][$j]
    }
}

$ws.Range("C2:O25").Value = $arr